# Actualización automática 2025-06-05 10:49:05
# Adds a new "CUMPLIMIENTO MENSUAL" worksheet (monthly compliance summary,
# grouped by product GRUPO instead of by CLIENTE) after the existing
# "VENTA MENSUAL" sheet.

$wb = $excel.ActiveWorkbook

$srcHeader = $wb.Worksheets.Item(1).Range("A1")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# ---- column widths ----
$ws.Columns.Item(1).ColumnWidth = 26
$ws.Columns.Item(2).ColumnWidth = 22
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 22
$ws.Columns.Item(6).ColumnWidth = 26

# ---- header row (reuse the bold/bordered header style already used by
# the other sheets, instead of re-building it from scratch) ----
$srcHeader.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

# ---- data rows, one per product group for this advisor ----
$asesor = "LINDAO ZUÑIGA BRYAN JOSE"

$grupos = @(
    @("240X120 PORCELANATO", 672.340305337043, 0, 672.340305337043, 0),
    @("240X80 PORCELANATO", 4992.1832, 0, 4992.1832, 0),
    @("FREGADEROS DE COCINA", 142.502095025027, 0, 142.502095025027, 0),
    @("GRANITO", 238.32, 0, 238.32, 0),
    @("GRIFERIAS", 106.82, 0, 106.82, 0),
    @("INODOROS", 2100, 0, 2100, 0),
    @("LAVABOS", 750, 0, 750, 0),
    @("LED", 300, 0, 300, 0),
    @("NO RESURTIBLES", 650.25, 9.58, 640.67, 0.01473279507881584),
    @("OTROS", 0, 0, 0, 0),
    @("PANELES DECORATIVOS", 350, 0, 350, 0),
    @("PANELES PU", 230, 0, 230, 0),
    @("PANELES PVC", 483, 0, 483, 0),
    @("PIEDRA SINTERIZADA", 1505.12, 0, 1505.12, 0),
    @("PORCELANATO", 38417.17, 233.38, 38183.79, 0.006074887869148092),
    @("PUERTAS DE SEGURIDAD", 342, 0, 342, 0),
    @("SAL SOLUBLE", 4130, 0, 4130, 0)
)

$row = 2
foreach ($g in $grupos) {
    $ws.Cells.Item($row, 1).Value = $asesor
    $ws.Cells.Item($row, 2).Value = $g[0]
    $ws.Cells.Item($row, 3).Value = $g[1]
    $ws.Cells.Item($row, 4).Value = $g[2]
    $ws.Cells.Item($row, 5).Value = $g[3]
    $ws.Cells.Item($row, 6).Value = $g[4]
    $row++
}

# ---- totals row ----
$totalRow = $row
$ws.Cells.Item($totalRow, 2).Value = "TOTAL"
$ws.Cells.Item($totalRow, 2).HorizontalAlignment = -4152
$ws.Cells.Item($totalRow, 3).Value = 55409.70560036207
$ws.Cells.Item($totalRow, 4).Value = 242.96
$ws.Cells.Item($totalRow, 5).Value = 55166.74560036208
$ws.Cells.Item($totalRow, 6).Value = 0.004384791389297914

# ---- number formats: currency for PRESUPUESTO/VENTA/POR CUMPLIR, percent
# for CUMPLIMIENTO (same numFmts already defined in the workbook) ----
$ws.Range("C2:E$totalRow").NumberFormat = "`"$`"#,##0.00"
$ws.Range("F2:F$totalRow").NumberFormat = "0.00%"

$ws.Range("A1").Select()
